# Natmi following Dr Hou advice
# Update LR-pair table rows 2-10 with refreshed NATMI statistics (3 clusters x 3 clusters)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Plau"
$ws.Cells.Item(2, 3).Value = "St14"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 43.21270066666667
$ws.Cells.Item(2, 8).Value = 129.638102
$ws.Cells.Item(2, 9).Value = 0.1487696778665633
$ws.Cells.Item(2, 10).Value = 0.1487696778665633
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.3701153333333333
$ws.Cells.Item(2, 14).Value = 1.110346
$ws.Cells.Item(2, 15).Value = 0.3716314658367922
$ws.Cells.Item(2, 16).Value = 0.3716314658367921
$ws.Cells.Item(2, 17).Value = 15.99368311147689
$ws.Cells.Item(2, 18).Value = 143.943148003292
$ws.Cells.Item(2, 19).Value = 0.0552874934576183
$ws.Cells.Item(2, 20).Value = 0.05528749345761828

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Plau"
$ws.Cells.Item(3, 3).Value = "St14"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 43.21270066666667
$ws.Cells.Item(3, 8).Value = 129.638102
$ws.Cells.Item(3, 9).Value = 0.1487696778665633
$ws.Cells.Item(3, 10).Value = 0.1487696778665633
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.5312866666666666
$ws.Cells.Item(3, 14).Value = 1.59386
$ws.Cells.Item(3, 15).Value = 0.5334630179589331
$ws.Cells.Item(3, 16).Value = 0.5334630179589331
$ws.Cells.Item(3, 17).Value = 22.95833169485778
$ws.Cells.Item(3, 18).Value = 206.62498525372
$ws.Cells.Item(3, 19).Value = 0.07936312133547516
$ws.Cells.Item(3, 20).Value = 0.07936312133547514

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Plau"
$ws.Cells.Item(4, 3).Value = "St14"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 43.21270066666667
$ws.Cells.Item(4, 8).Value = 129.638102
$ws.Cells.Item(4, 9).Value = 0.1487696778665633
$ws.Cells.Item(4, 10).Value = 0.1487696778665633
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.09451833333333333
$ws.Cells.Item(4, 14).Value = 0.283555
$ws.Cells.Item(4, 15).Value = 0.09490551620427472
$ws.Cells.Item(4, 16).Value = 0.0949055162042747
$ws.Cells.Item(4, 17).Value = 4.084392445845555
$ws.Cells.Item(4, 18).Value = 36.75953201261
$ws.Cells.Item(4, 19).Value = 0.01411906307346985
$ws.Cells.Item(4, 20).Value = 0.01411906307346985

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Plau"
$ws.Cells.Item(5, 3).Value = "St14"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 237.7114053333333
$ws.Cells.Item(5, 8).Value = 713.134216
$ws.Cells.Item(5, 9).Value = 0.8183762794517323
$ws.Cells.Item(5, 10).Value = 0.8183762794517323
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.3701153333333333
$ws.Cells.Item(5, 14).Value = 1.110346
$ws.Cells.Item(5, 15).Value = 0.3716314658367922
$ws.Cells.Item(5, 16).Value = 0.3716314658367921
$ws.Cells.Item(5, 17).Value = 87.98063602208177
$ws.Cells.Item(5, 18).Value = 791.8257241987359
$ws.Cells.Item(5, 19).Value = 0.3041343763387075
$ws.Cells.Item(5, 20).Value = 0.3041343763387075

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Plau"
$ws.Cells.Item(6, 3).Value = "St14"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 237.7114053333333
$ws.Cells.Item(6, 8).Value = 713.134216
$ws.Cells.Item(6, 9).Value = 0.8183762794517323
$ws.Cells.Item(6, 10).Value = 0.8183762794517323
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.5312866666666666
$ws.Cells.Item(6, 14).Value = 1.59386
$ws.Cells.Item(6, 15).Value = 0.5334630179589331
$ws.Cells.Item(6, 16).Value = 0.5334630179589331
$ws.Cells.Item(6, 17).Value = 126.2929001681955
$ws.Cells.Item(6, 18).Value = 1136.63610151376
$ws.Cells.Item(6, 19).Value = 0.4365734798623244
$ws.Cells.Item(6, 20).Value = 0.4365734798623244

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Plau"
$ws.Cells.Item(7, 3).Value = "St14"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 237.7114053333333
$ws.Cells.Item(7, 8).Value = 713.134216
$ws.Cells.Item(7, 9).Value = 0.8183762794517323
$ws.Cells.Item(7, 10).Value = 0.8183762794517323
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.09451833333333333
$ws.Cells.Item(7, 14).Value = 0.283555
$ws.Cells.Item(7, 15).Value = 0.09490551620427472
$ws.Cells.Item(7, 16).Value = 0.0949055162042747
$ws.Cells.Item(7, 17).Value = 22.46808584643111
$ws.Cells.Item(7, 18).Value = 202.21277261788
$ws.Cells.Item(7, 19).Value = 0.07766842325070043
$ws.Cells.Item(7, 20).Value = 0.07766842325070042

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Plau"
$ws.Cells.Item(8, 3).Value = "St14"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 9.543019333333334
$ws.Cells.Item(8, 8).Value = 28.629058
$ws.Cells.Item(8, 9).Value = 0.03285404268170446
$ws.Cells.Item(8, 10).Value = 0.03285404268170446
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.3701153333333333
$ws.Cells.Item(8, 14).Value = 1.110346
$ws.Cells.Item(8, 15).Value = 0.3716314658367922
$ws.Cells.Item(8, 16).Value = 0.3716314658367921
$ws.Cells.Item(8, 17).Value = 3.532017781563111
$ws.Cells.Item(8, 18).Value = 31.788160034068
$ws.Cells.Item(8, 19).Value = 0.01220959604046636
$ws.Cells.Item(8, 20).Value = 0.01220959604046636

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Plau"
$ws.Cells.Item(9, 3).Value = "St14"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 9.543019333333334
$ws.Cells.Item(9, 8).Value = 28.629058
$ws.Cells.Item(9, 9).Value = 0.03285404268170446
$ws.Cells.Item(9, 10).Value = 0.03285404268170446
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.5312866666666666
$ws.Cells.Item(9, 14).Value = 1.59386
$ws.Cells.Item(9, 15).Value = 0.5334630179589331
$ws.Cells.Item(9, 16).Value = 0.5334630179589331
$ws.Cells.Item(9, 17).Value = 5.070078931542222
$ws.Cells.Item(9, 18).Value = 45.63071038387999
$ws.Cells.Item(9, 19).Value = 0.01752641676113366
$ws.Cells.Item(9, 20).Value = 0.01752641676113366

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Plau"
$ws.Cells.Item(10, 3).Value = "St14"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 9.543019333333334
$ws.Cells.Item(10, 8).Value = 28.629058
$ws.Cells.Item(10, 9).Value = 0.03285404268170446
$ws.Cells.Item(10, 10).Value = 0.03285404268170446
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.09451833333333333
$ws.Cells.Item(10, 14).Value = 0.283555
$ws.Cells.Item(10, 15).Value = 0.09490551620427472
$ws.Cells.Item(10, 16).Value = 0.0949055162042747
$ws.Cells.Item(10, 17).Value = 0.9019902823544445
$ws.Cells.Item(10, 18).Value = 8.11791254119
$ws.Cells.Item(10, 19).Value = 0.003118029880104436
$ws.Cells.Item(10, 20).Value = 0.003118029880104436

